$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.697.08"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "3.324.79"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.06"
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.84"
$ws.Range("E6").Value = "  -7.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -2.79%  "
$ws.Range("D9").Value = "3.321.17"
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("E10").Value = "  -5.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.575"
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.44"
$ws.Range("E12").Value = "  -4.63%  "
$ws.Range("E13").Value = "  -4.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "669.72"
$ws.Range("E14").Value = "  +4.75%  "
$ws.Range("D15").Value = "3.859.75"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.38"
$ws.Range("E16").Value = "  -3.11%  "
$ws.Range("D17").Value = "67.851.78"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").Value = "3.331.13"
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("E20").Value = "  -3.90%  "
$ws.Range("E21").Value = "  -2.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.888"
$ws.Range("E22").Value = "  -2.83%  "
$ws.Range("E23").Value = "  +5.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.97"
$ws.Range("E24").Value = "  -5.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.67"
$ws.Range("E25").Value = "  -2.24%  "
$ws.Range("E26").Value = "  -5.13%  "
$ws.Range("E27").Value = "  -6.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.30"
$ws.Range("E28").Value = "  -4.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.41"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("E30").Value = "  -4.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.33"
$ws.Range("E31").Value = "  +6.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "585.40"
$ws.Range("E32").Value = "  -4.33%  "
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("E34").Value = "  -2.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.82"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.26"
$ws.Range("E38").Value = "  -14.82%  "
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "32.78"
$ws.Range("E40").Value = "  -3.28%  "
$ws.Range("E42").Value = "  -5.34%  "
$ws.Range("E43").Value = "  -3.69%  "
$ws.Range("E44").Value = "  -6.29%  "
$ws.Range("E45").Value = "  -5.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0407"
$ws.Range("E46").Value = "  -3.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.59"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  -4.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "126.87"
$ws.Range("E51").Value = "  -0.84%  "

# Row 35/36: content swapped between Dai and Maker, with updated price/volume data
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "3.717.47"
$ws.Range("E35").Value = "  -8.36%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.24%  "
